$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, border, centered) from the last existing
# header cell (AC1) into the three new header cells, then set their text.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))

$ws.Range("AD1").Value() = "Wins"
$ws.Range("AE1").Value() = "Losses"
$ws.Range("AF1").Value() = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row.
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value() = 90
    $ws.Cells.Item($r, 31).Value() = 72
    $ws.Cells.Item($r, 32).Value() = 0
}
